$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "ItemType" header (W1) to "DenominationName"
$ws.Range("W1").Value = "DenominationName"

# Insert a new column at X1 (shifting Count/Amount one column to the right)
# and give it the new "DenominationValue" header.
$ws.Range("X1").Insert(-4161) | Out-Null
$ws.Range("X1").Value = "DenominationValue"
